# Apply updated cryptocurrency price/volume data (GitHub Actions scheduled refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For values that look like plain numbers (e.g. "599.69"), prefix with a literal
# apostrophe so Excel stores them as text, matching the other price cells in column D
# which already contain text such as "65.376.71" (thousands-dot formatted).

$ws.Range('D2').Value = '65.376.71'
$ws.Range('E2').Value = '  +2.60%  '

$ws.Range('D3').Value = '3.203.07'
$ws.Range('E3').Value = '  +1.84%  '

$ws.Range('E4').Value = '  -0.13%  '

$ws.Range('D5').Value = '''599.69'
$ws.Range('E5').Value = '  +2.26%  '

$ws.Range('D6').Value = '''153.88'
$ws.Range('E6').Value = '  +5.94%  '

$ws.Range('E7').Value = '  -0.06%  '

$ws.Range('D8').Value = '3.201.92'
$ws.Range('E8').Value = '  +2.06%  '

$ws.Range('E9').Value = '  +2.86%  '

$ws.Range('D10').Value = '''0.168'
$ws.Range('E10').Value = '  +4.61%  '

$ws.Range('E11').Value = '  +6.48%  '

$ws.Range('D12').Value = '''0.473'
$ws.Range('E12').Value = '  +3.07%  '

$ws.Range('D13').Value = '''0.0000256'
$ws.Range('E13').Value = '  +3.87%  '

$ws.Range('D14').Value = '''39.39'
$ws.Range('E14').Value = '  +6.84%  '

$ws.Range('D15').Value = '3.731.19'
$ws.Range('E15').Value = '  +1.56%  '

$ws.Range('E16').Value = '  +0.47%  '

$ws.Range('D17').Value = '''7.41'
$ws.Range('E17').Value = '  +4.39%  '

$ws.Range('D18').Value = '65.040.09'
$ws.Range('E18').Value = '  +2.37%  '

$ws.Range('D19').Value = '3.201.42'
$ws.Range('E19').Value = '  +1.58%  '

$ws.Range('D20').Value = '''485.43'
$ws.Range('E20').Value = '  +4.73%  '

$ws.Range('D21').Value = '''15.06'
$ws.Range('E21').Value = '  +5.68%  '

$ws.Range('D22').Value = '''0.774'
$ws.Range('E22').Value = '  +6.17%  '

$ws.Range('D23').Value = '''7.92'
$ws.Range('E23').Value = '  +6.41%  '

$ws.Range('B24').Value = 'InternetComputer(DFINITY)'
$ws.Range('C24').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D24').Value = '''13.85'
$ws.Range('E24').Value = '  +6.61%  '

$ws.Range('B25').Value = 'Fetch.AI'
$ws.Range('C25').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D25').Value = '''2.45'
$ws.Range('E25').Value = '  +11.71%  '

$ws.Range('D26').Value = '''83.54'
$ws.Range('E26').Value = '  +2.73%  '

$ws.Range('E27').Value = '  +0.41%  '

$ws.Range('D28').Value = '''9.81'
$ws.Range('E28').Value = '  +7.63%  '

$ws.Range('E29').Value = '  +4.01%  '

$ws.Range('E30').Value = '  +4.58%  '

$ws.Range('D31').Value = '''7.50'
$ws.Range('E31').Value = '  +7.38%  '

$ws.Range('E32').Value = '  -0.01%  '

$ws.Range('D33').Value = '''0.121'
$ws.Range('E33').Value = '  +9.34%  '

$ws.Range('D34').Value = '''28.64'
$ws.Range('E34').Value = '  +6.36%  '

$ws.Range('D35').Value = '0.0₃0903'
$ws.Range('E35').Value = '  +5.60%  '

$ws.Range('D36').Value = '''3.65'
$ws.Range('E36').Value = '  +8.28%  '

$ws.Range('E37').Value = '  +4.44%  '

$ws.Range('B38').Value = 'Filecoin'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D38').Value = '''6.36'
$ws.Range('E38').Value = '  +5.88%  '

$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').Value = '''2.38'
$ws.Range('E39').Value = '  +3.55%  '

$ws.Range('D40').Value = '''480.90'
$ws.Range('E40').Value = '  +9.37%  '

$ws.Range('D41').Value = '''52.20'
$ws.Range('E41').Value = '  +3.18%  '

$ws.Range('D42').Value = '''9.45'
$ws.Range('E42').Value = '  +8.38%  '

$ws.Range('E43').Value = '  +10.29%  '

$ws.Range('D44').Value = '''0.0384'
$ws.Range('E44').Value = '  +3.56%  '

$ws.Range('D45').Value = '2.952.37'
$ws.Range('E45').Value = '  +1.46%  '

$ws.Range('E46').Value = '  +5.12%  '

$ws.Range('D47').Value = '''38.91'
$ws.Range('E47').Value = '  +7.52%  '

$ws.Range('D48').Value = '''131.87'
$ws.Range('E48').Value = '  +5.03%  '

$ws.Range('E49').Value = '  +7.67%  '

$ws.Range('D50').Value = '''25.67'
$ws.Range('E50').Value = '  +5.70%  '

$ws.Range('E51').Value = '  -0.02%  '
